# Update cryptocurrency price/volume(1h) figures and reorder two coin pairs
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.150.43'
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("D3").Value = '2.209.71'
$ws.Range("E3").Value = '  -0.71%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '295.85'
$ws.Range("E5").Value = '  +1.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.75'
$ws.Range("E6").Value = '  +0.36%  '
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  -0.63%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.19'
$ws.Range("E10").Value = '  +7.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '30.94'
$ws.Range("E11").Value = '  +1.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0781'
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("E13").Value = '  +2.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.39'
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").Value = '2.550.96'
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.81'
$ws.Range("E16").Value = '  -0.94%  '
$ws.Range("D17").Value = '2.236.44'
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("D19").Value = '40.063.75'
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").Value = '0.0₃0887'
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.32'
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.76'
$ws.Range("E22").Value = '  -1.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.66'
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.37'
$ws.Range("E24").Value = '  -0.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  +1.03%  '
$ws.Range("E27").Value = '  -0.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.19'
$ws.Range("E28").Value = '  +2.26%  '
$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.32'
$ws.Range("E29").Value = '  +1.15%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.08'
$ws.Range("E30").Value = '  -4.96%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '156.27'
$ws.Range("E31").Value = '  -0.07%  '
$ws.Range("E32").Value = '  +0.99%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.97'
$ws.Range("E34").Value = '  +0.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.00'
$ws.Range("E35").Value = '  +3.92%  '
$ws.Range("E36").Value = '  -0.57%  '
$ws.Range("E37").Value = '  -0.83%  '
$ws.Range("E38").Value = '  +1.55%  '
$ws.Range("E39").Value = '  +3.01%  '
$ws.Range("E40").Value = '  +2.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.60'
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.81'
$ws.Range("E42").Value = '  -1.30%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '2.070.73'
$ws.Range("E43").Value = '  -2.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.14'
$ws.Range("E44").Value = '  +5.42%  '
$ws.Range("E45").Value = '  +1.06%  '
$ws.Range("E46").Value = '  +0.27%  '
$ws.Range("E47").Value = '  +5.58%  '
$ws.Range("E48").Value = '  -11.13%  '
$ws.Range("D49").Value = '2.424.34'
$ws.Range("E49").Value = '  -0.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.13'
$ws.Range("E50").Value = '  +2.37%  '
$ws.Range("E51").Value = '  +0.92%  '
